$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad"
$ws.Range("G4").Value = "Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud"
$ws.Range("G5").Value = "Dr. Menna tu'Alllah Mohammad, Dr. Abeer Ragab"
$ws.Range("G9").Value = "Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Eman M. Abo-Sakaya"
$ws.Range("G13").Value = "Dr. Youstina Gamil, Dr. Sarah Mahdy"
$ws.Range("G14").Value = "Dr. Alaa Ashraf, Dr. Sarah Mahdy"
$ws.Range("G15").Value = "Dr. Alaa Ashraf, Dr. Sarah Mahdy"
$ws.Range("G18").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida"
$ws.Range("G19").Value = "Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Sorial, Dr. Wafaa Ebida"
$ws.Range("G20").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G21").Value = "Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad"
$ws.Range("G23").Value = "Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Heba Mahmoud Ali, Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud"
$ws.Range("G24").Value = "Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Lamiaa Ossama, Dr. Abeer Ragab, Dr. Amera Ahmad Saad"
$ws.Range("G25").Value = "Dr. Yasmin Tarek, Dr. Nourhan Mohammad"
$ws.Range("G28").Value = "Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Marwa Mustafa, Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Esraa Mostafa, Dr. Nourhan Osama"
$ws.Range("G29").Value = "Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Esraa Mostafa"
$ws.Range("G34").Value = "Dr. Alaa Ashraf, Dr. Sarah Mahdy"
$ws.Range("G37").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida"
$ws.Range("G38").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G39").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G40").Value = "Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Hend Mahmoud, Dr. Heba Mahmoud Ali, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad"
$ws.Range("G41").Value = "Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Amira Sobhy"
$ws.Range("G42").Value = "Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad"
$ws.Range("G43").Value = "Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Lamiaa Ossama, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G47").Value = "Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Esraa Mostafa, Dr. Nourhan Osama"
$ws.Range("G48").Value = "Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad"
$ws.Range("G49").Value = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Range("G52").Value = "Dr. Alaa Ashraf, Dr. Sarah Mahdy"
$ws.Range("G54").Value = "Dr. Afaf Abdallah, Dr. Amr Saeed"
$ws.Range("G56").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida"
$ws.Range("G57").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G58").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G59").Value = "Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Servinaz Sayed Mohammad"
$ws.Range("G60").Value = "Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Amira Sobhy"
$ws.Range("G63").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G66").Value = "Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amira Ibrahim, Dr. Marina Youhanna, Dr. Madeha Saeed"
$ws.Range("G67").Value = "Dr. Yasmeena Fattoh, Dr. Amira Ibrahim, Dr. Esraa Mostafa"
$ws.Range("G75").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida"
$ws.Range("G76").Value = "Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Sorial, Dr. Wafaa Ebida"
$ws.Range("G77").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G78").Value = "Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Servinaz Sayed Mohammad"
$ws.Range("G79").Value = "Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Amira Sobhy"
$ws.Range("G81").Value = "Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Lamiaa Ossama, Dr. Abeer Ragab, Dr. Amera Ahmad Saad"
$ws.Range("G82").Value = "Dr. Yasmin Tarek, Dr. Nourhan Mohammad"
$ws.Range("G83").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G85").Value = "Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Eman M. Abo-Sakaya, Dr. Amira Ibrahim, Dr. Marina Youhanna, Dr. Madeha Saeed"
$ws.Range("G86").Value = "Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Amany Raafat, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Maryam Ahmad"
$ws.Range("G89").Value = "Dr. Youstina Gamil, Dr. Sarah Mahdy"
$ws.Range("G94").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida"
$ws.Range("G95").Value = "Dr. Neveen Nashaat, Dr. Yasmin, Dr. Eman Samir Gabry, Dr. Marina Sorial, Dr. Wafaa Ebida"
$ws.Range("G96").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
$ws.Range("G97").Value = "Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Heba Mahmoud Ali, Dr. Amira Sobhy, Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Servinaz Sayed Mohammad"
$ws.Range("G98").Value = "Dr. Alshimaa Atef, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Amira Sobhy"
$ws.Range("G99").Value = "Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Servinaz Sayed Mohammad"
$ws.Range("G100").Value = "Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Kerelos Zareef, Dr. Lamiaa Ossama, Dr. Abeer Ragab, Dr. Amera Ahmad Saad, Dr. Menna tu'Alllah Mohammad"
$ws.Range("G101").Value = "Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed"
$ws.Range("G104").Value = "Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim, Dr. Maryam Ahmad, Dr. Esraa Mostafa, Dr. Nourhan Osama"
$ws.Range("G113").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida"
$ws.Range("G115").Value = "Dr. Neveen Nashaat, Dr. Nardine, Dr. Yasmin, Dr. Remon, Dr. Monica, Dr. Marina Sorial, Dr. Marina Atef"
